# Add four new conversation rows (13-16) to the conversation log sheet,
# matching the data added by the commit (new inbound/outbound text messages).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Timestamp, Sender, Sender Id, Phone, Message
$newRows = @(
    @("2025-09-21 14:24:59", "Noah", 8450689526, "13052054965", "Hello"),
    @("2025-09-21 15:40:12", "Noah", 8450689526, "13052054965", "Hello"),
    @("2025-09-21 16:03:15", "Noah", 8450689526, "13052054965", "Hello"),
    @("2025-09-21 16:04:02", "Noah", 8450689526, "13052054965", "Hey man")
)

$startRow = 13

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $entry = $newRows[$i]

    # Timestamp (text)
    $ws.Cells.Item($r, 1).Value = $entry[0]

    # Sender (text)
    $ws.Cells.Item($r, 2).Value = $entry[1]

    # Sender Id (numeric)
    $ws.Cells.Item($r, 3).Value = $entry[2]

    # Phone - keep as text even though it is numeric-looking
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $entry[3]

    # Message (text)
    $ws.Cells.Item($r, 5).Value = $entry[4]

    # Media / Channel columns stay empty (blank inline strings), same as
    # every other row in the sheet. A leading "'" forces the cell to be
    # recorded as an (empty) text value instead of being dropped entirely,
    # then the quote-prefix formatting is cleared so the cell keeps the
    # default style used by the rest of the table.
    $ws.Cells.Item($r, 6).Value = "'"
    $ws.Cells.Item($r, 6).Style = "Normal"
    $ws.Cells.Item($r, 7).Value = "'"
    $ws.Cells.Item($r, 7).Style = "Normal"
}
